$d = $word.ActiveDocument

# 1. Remove the "Meta description" paragraph that immediately follows the
#    title heading (it is being dropped entirely from the top of the doc).
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Meta description:*") {
        $null = $p.Range.Delete()
        $found = $true
        break
    }
}

# 2. The closing "Prompt: ..." paragraph is replaced by two paragraphs: a
#    bold restatement of the page title, followed by the (now-relocated)
#    meta-description text in italics.
$count = $d.Paragraphs.Count
$pLast = $d.Paragraphs.Item($count)
$null = $pLast.Range.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play 100 Fortunes Free and Enjoy Expansive Chinese-Themed Slots</w:t></w:r></w:p><w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of 100 Fortunes, an expanding slot game with a unique Chinese aesthetic. Play for free and enjoy bonus features like free spins and expanding Wild symbols.</w:t></w:r></w:p>")

Write-Output "metaDescriptionRemoved=$found"
Write-Output "finalParagraphCount=$($d.Paragraphs.Count)"
